$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The surviving data row (former row 2) keeps its Target cluster label as a
# raw shared-string index; once the "ECs" string entry is dropped below,
# that same index resolves to "MuSCs" instead. Setting the cell text here
# keeps behaviour correct regardless of how the engine renumbers strings.
$ws.Range("D2").Value = "MuSCs"

# Recomputed (new TPM) values for the surviving row.
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02430333333333333
$ws.Range("N2").Value = 0.07291
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.00016521406
$ws.Range("R2").Value = 0.00148692654
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove the now-obsolete third row (the ECs/MuSCs pair collapses to a
# single surviving row) so the sheet dimension shrinks back to A1:T2.
$ws.Rows.Item(3).Delete()
